$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 27.7608239223043
$ws.Cells.Item(2, 3).Value = 30.43167762435635
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 27.76415565283329
$ws.Cells.Item(3, 3).Value = 30.43427447041097
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 27.80085108117647
$ws.Cells.Item(4, 3).Value = 30.38770512775293
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 27.76251901515253
$ws.Cells.Item(5, 3).Value = 30.4021683580043
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 27.76072746189019
$ws.Cells.Item(6, 3).Value = 30.42522454420415
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 27.77727973682701
$ws.Cells.Item(7, 3).Value = 30.41906045070446
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 27.78446173692694
$ws.Cells.Item(8, 3).Value = 30.41462334204979
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 27.76050118429297
$ws.Cells.Item(9, 3).Value = 30.41926158129405
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 27.76686598914338
$ws.Cells.Item(10, 3).Value = 30.41849637698485
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 27.77808232087939
$ws.Cells.Item(11, 3).Value = 30.43074601113019
$ws.Cells.Item(12, 2).Value = 27.77162681014265
$ws.Cells.Item(12, 3).Value = 30.4183237886892

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 19.71096378066319
$ws.Cells.Item(2, 3).Value = 27.12554861644868
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 19.70666034781575
$ws.Cells.Item(3, 3).Value = 27.12441247315927
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 19.69473720740051
$ws.Cells.Item(4, 3).Value = 27.10100171302534
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 19.70634436507044
$ws.Cells.Item(5, 3).Value = 27.14092747501897
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 19.71380722760305
$ws.Cells.Item(6, 3).Value = 27.10969116205595
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 19.6990256308485
$ws.Cells.Item(7, 3).Value = 27.14155728141581
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 19.69283585213583
$ws.Cells.Item(8, 3).Value = 27.11998882156776
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 19.7109287530526
$ws.Cells.Item(9, 3).Value = 27.14285905819212
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 19.70808572588459
$ws.Cells.Item(10, 3).Value = 27.1488540439452
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 19.69632878342746
$ws.Cells.Item(11, 3).Value = 27.13944556714049
$ws.Cells.Item(12, 2).Value = 19.70397176739019
$ws.Cells.Item(12, 3).Value = 27.12942862119696

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 14.71210866387949
$ws.Cells.Item(2, 3).Value = 22.85673989991717
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 14.71457888402703
$ws.Cells.Item(3, 3).Value = 22.80292029071429
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 14.69337827589535
$ws.Cells.Item(4, 3).Value = 22.78669765343427
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 14.70796967197167
$ws.Cells.Item(5, 3).Value = 22.82267012730443
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 14.70728687960158
$ws.Cells.Item(6, 3).Value = 22.84514864959789
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 14.68214970443914
$ws.Cells.Item(7, 3).Value = 22.85584582650917
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 14.69367101778813
$ws.Cells.Item(8, 3).Value = 22.81512895830863
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 14.69523298733144
$ws.Cells.Item(9, 3).Value = 22.88361392675698
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 14.68916523002296
$ws.Cells.Item(10, 3).Value = 22.8043320559091
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 14.69387074554603
$ws.Cells.Item(11, 3).Value = 22.81169227272106
$ws.Cells.Item(12, 2).Value = 14.69894120605028
$ws.Cells.Item(12, 3).Value = 22.8284789661173

Write-Output "Updated sigma_010, sigma_025, sigma_050 sheets"